$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: A5 stays "Dall Makhini"; C5's messy Google image-redirect
#     URL is replaced below with a clean, hyperlink-styled image URL. ---
$ws.Range("A5").Value = "Dall Makhini"

# --- Row 6: brand new dish entry appended below the existing data. ---
$ws.Range("C6").Value = "https://img-global.cpcdn.com/recipes/7650b1942bfd97c4/400x400cq70/photo.jpg"
$ws.Range("A6").Value = "butter naan dal makhani"

# C5 gets the new image URL, styled with Excel's built-in "Hyperlink"
# cell style (underlined, theme-coloured font) -- added then immediately
# unlinked so the cell keeps the visual style without a live hyperlink.
$ws.Hyperlinks.Add($ws.Range("C5"), "https://shwetainthekitchen.com/wp-content/uploads/2019/11/IMG_6917_1-scaled.jpg", "", "", "https://shwetainthekitchen.com/wp-content/uploads/2019/11/IMG_6917_1-scaled.jpg")
$ws.Hyperlinks.Delete()

$ws.Range("B6").Value = 180

# Leave the selection where the author last left it.
$ws.Range("F18").Select() | Out-Null
